# Apply scheduled market-data refresh updates to the crafting profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1709.6389
$ws.Range("I40").Value = 1461.8096
$ws.Range("J40").Value = 2056.6
$ws.Range("K40").Value = 1461.8096
$ws.Range("L40").Value = 2056.6
$ws.Range("M40").Value = -1286.8096
$ws.Range("N40").Value = -2406.6
$ws.Range("H44").Value = 30000
$ws.Range("J44").Value = 30000
$ws.Range("L44").Value = 30000
$ws.Range("N44").Value = -30924
$ws.Range("H64").Value = 3812.5557
$ws.Range("I64").Value = 3662.6667
$ws.Range("J64").Value = 3999.9167
$ws.Range("K64").Value = 3662.6667
$ws.Range("L64").Value = 3999.9167
$ws.Range("M64").Value = -3414.6667
$ws.Range("N64").Value = -4495.9167
$ws.Range("H67").Value = 3812.5557
$ws.Range("I67").Value = 3662.6667
$ws.Range("J67").Value = 3999.9167
$ws.Range("K67").Value = 3662.6667
$ws.Range("L67").Value = 3999.9167
$ws.Range("M67").Value = -2804.6667
$ws.Range("N67").Value = -5715.9167
$ws.Range("H70").Value = 1547.5
$ws.Range("I70").Value = 1078
$ws.Range("J70").Value = 1704
$ws.Range("K70").Value = 3234
$ws.Range("L70").Value = 5112
$ws.Range("M70").Value = -2964
$ws.Range("N70").Value = -5652
$ws.Range("H73").Value = 1547.5
$ws.Range("I73").Value = 1078
$ws.Range("J73").Value = 1704
$ws.Range("K73").Value = 3234
$ws.Range("L73").Value = 5112
$ws.Range("M73").Value = -2298
$ws.Range("N73").Value = -6984
$ws.Range("H76").Value = 3020
$ws.Range("I76").Value = 3020.513
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 3020.513
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -2705.513
$ws.Range("N76").Value = -3630
$ws.Range("H79").Value = 3020
$ws.Range("I79").Value = 3020.513
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 3020.513
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -1928.513
$ws.Range("N79").Value = -5184
$ws.Range("H82").Value = 994.1111
$ws.Range("I82").Value = 994.1111
$ws.Range("K82").Value = 2982.3333
$ws.Range("M82").Value = -2576.3333
$ws.Range("H85").Value = 994.1111
$ws.Range("I85").Value = 994.1111
$ws.Range("K85").Value = 2982.3333
$ws.Range("M85").Value = -1578.3333
$ws.Range("H88").Value = 3396.2917
$ws.Range("I88").Value = 1854.4286
$ws.Range("J88").Value = 4031.1765
$ws.Range("K88").Value = 1854.4286
$ws.Range("L88").Value = 4031.1765
$ws.Range("M88").Value = -1448.4286
$ws.Range("N88").Value = -4843.1765
$ws.Range("H91").Value = 3396.2917
$ws.Range("I91").Value = 1854.4286
$ws.Range("J91").Value = 4031.1765
$ws.Range("K91").Value = 1854.4286
$ws.Range("L91").Value = 4031.1765
$ws.Range("M91").Value = -450.4286
$ws.Range("N91").Value = -6839.1765
$ws.Range("H97").Value = 40758.2
$ws.Range("J97").Value = 40758.2
$ws.Range("L97").Value = 122274.6
$ws.Range("N97").Value = -123266.6
$ws.Range("H116").Value = 4126.769
$ws.Range("I116").Value = 2227.4375
$ws.Range("J116").Value = 7165.7
$ws.Range("K116").Value = 2227.4375
$ws.Range("L116").Value = 7165.7
$ws.Range("M116").Value = 1214.5625
$ws.Range("N116").Value = -14049.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1643.5834
$ws.Range("I45").Value = 1491.6842
$ws.Range("K45").Value = 1491.6842
$ws.Range("M45").Value = -1114.6842
$ws.Range("H74").Value = 86098.08
$ws.Range("I74").Value = 112424.555
$ws.Range("J74").Value = 26863.5
$ws.Range("K74").Value = 112424.555
$ws.Range("L74").Value = 26863.5
$ws.Range("M74").Value = -111550.555
$ws.Range("N74").Value = -28611.5
$ws.Range("H77").Value = 86098.08
$ws.Range("I77").Value = 112424.555
$ws.Range("J77").Value = 26863.5
$ws.Range("K77").Value = 562122.7749999999
$ws.Range("L77").Value = 134317.5
$ws.Range("M77").Value = -557754.7749999999
$ws.Range("N77").Value = -143053.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3512
$ws.Range("I134").Value = 3145.756
$ws.Range("J134").Value = 5657.143
$ws.Range("K134").Value = 9437.268
$ws.Range("L134").Value = 16971.429
$ws.Range("M134").Value = -6902.268
$ws.Range("N134").Value = -22041.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1013.95654
$ws.Range("I122").Value = 537.9091
$ws.Range("J122").Value = 1450.3334
$ws.Range("K122").Value = 4841.1819
$ws.Range("L122").Value = 13053.0006
$ws.Range("M122").Value = -2391.1819
$ws.Range("N122").Value = -17953.0006
$ws.Range("H131").Value = 884.32
$ws.Range("I131").Value = 526
$ws.Range("J131").Value = 903.17896
$ws.Range("K131").Value = 1578
$ws.Range("L131").Value = 2709.53688
$ws.Range("M131").Value = 3462
$ws.Range("N131").Value = -12789.53688
$ws.Range("H132").Value = 2665.9546
$ws.Range("I132").Value = 2603.1177
$ws.Range("J132").Value = 2879.6
$ws.Range("K132").Value = 23428.0593
$ws.Range("L132").Value = 25916.4
$ws.Range("M132").Value = -20898.0593
$ws.Range("N132").Value = -30976.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2693.5667
$ws.Range("I132").Value = 2327.389
$ws.Range("J132").Value = 3242.8333
$ws.Range("K132").Value = 6982.167
$ws.Range("L132").Value = 9728.499899999999
$ws.Range("M132").Value = -4452.167
$ws.Range("N132").Value = -14788.4999
$ws.Range("H136").Value = 2205.8462
$ws.Range("I136").Value = 2016
$ws.Range("J136").Value = 3250
$ws.Range("K136").Value = 6048
$ws.Range("L136").Value = 9750
$ws.Range("M136").Value = -3498
$ws.Range("N136").Value = -14850

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1487.5
$ws.Range("I81").Value = 1185
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 2370
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -1309
$ws.Range("N81").Value = -8122
$ws.Range("H84").Value = 1487.5
$ws.Range("I84").Value = 1185
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 11850
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -6546
$ws.Range("N84").Value = -40608
$ws.Range("H136").Value = 4853.294
$ws.Range("I136").Value = 5955.1816
$ws.Range("J136").Value = 2833.1667
$ws.Range("K136").Value = 17865.5448
$ws.Range("L136").Value = 8499.500100000001
$ws.Range("M136").Value = -15315.5448
$ws.Range("N136").Value = -13599.5001
